$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Duplicate slide 6 ("What is Multi-Threading?") to create the new slide 7
#    ("Overiew") while it still carries the original layout/formatting
#    (title normAutofit, grpSpPr xfrm, extLst creationId, etc.)
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Copy()
$p.Slides.Paste(7) | Out-Null
$s7 = $p.Slides.Item(7)

# ---------------------------------------------------------------------------
# 2. New slide 7 content: "Overiew" / "Lots of bad jokes (I can’t help it)"
# ---------------------------------------------------------------------------
$t7 = $s7.Shapes.Item(1).TextFrame.TextRange
$t7.Text = "Overiew"

$c7 = $s7.Shapes.Item(2).TextFrame.TextRange
$c7.Text = "Lots of bad jokes (I can’t help it)"
$c7.Characters(1, 5).Text = "Lots "
$c7.InsertAfter("`r") | Out-Null

# ---------------------------------------------------------------------------
# 3. Slide 6 title: "Broism: Bilbro Baggins" (4 runs)
# ---------------------------------------------------------------------------
$title6 = $s6.Shapes.Item(1).TextFrame.TextRange
$title6.Text = "Broism: Bilbro Baggins"
# Broism(1-6) ": "(7-8) Bilbro(9-14) " Baggins"(15-22)
$title6.Characters(1, 6).Text = "Broism"
$title6.Characters(7, 2).Text = ": "
$title6.Characters(9, 6).Text = "Bilbro"
$title6.Characters(15, 8).Text = " Baggins"

# ---------------------------------------------------------------------------
# 4. Slide 6 content placeholder: "Bilbro Baggins" definition paragraph
# ---------------------------------------------------------------------------
$body6 = $s6.Shapes.Item(2).TextFrame.TextRange
$desc = ":  Your bro who is obsessed with Lord of the Rings.  Example:  Joe has seen the twin towers like 5 times.  He’s such a "
$full = "Bilbro" + " " + "Baggins" + $desc + "Bilbro" + " Baggins" + "."
$body6.Text = $full

$pos = 1
$r1 = $body6.Characters($pos, 6); $r1.Font.Bold = $true; $pos += 6          # Bilbro
$r2 = $body6.Characters($pos, 1); $r2.Font.Bold = $true; $pos += 1         # " "
$r3 = $body6.Characters($pos, 7); $r3.Font.Bold = $true; $pos += 7         # Baggins
$r4 = $body6.Characters($pos, $desc.Length); $pos += $desc.Length          # description
$r5 = $body6.Characters($pos, 6); $r5.Font.Italic = $true; $pos += 6       # Bilbro
$r6 = $body6.Characters($pos, 8); $r6.Font.Italic = $true; $pos += 8       # " Baggins"
$r7 = $body6.Characters($pos, 1); $pos += 1                                # "."

# Remove the bullet and flush the paragraph's indent to the left margin
# (marL="0" indent="0" + buNone), matching a manual "No bullet" toggle.
$body6.ParagraphFormat.Bullet.Visible = $false
$ruler6 = $s6.Shapes.Item(2).TextFrame.Ruler
$lvl6 = $ruler6.Levels.Item(1)
$lvl6.LeftMargin = 0
$lvl6.FirstMargin = 0
